$d = $word.ActiveDocument

# 1. Remove the first paragraph ("ĐỀ:") entirely, including its paragraph mark.
$p1 = $d.Paragraphs(1).Range
$p1.Delete()

Write-Output "done"
